$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-25 Friday" "2025-07-26 Saturday"

Replace-Text "242×8=" "924×2="
Replace-Text "455×7=" "184×4="
Replace-Text "121×9=" "128×7="
Replace-Text "303×2=" "251×9="
Replace-Text "808×2=" "443×3="
Replace-Text "567×2=" "758×4="
Replace-Text "870×9=" "304×2="
Replace-Text "346×2=" "311×6="
Replace-Text "664×9=" "879×9="
Replace-Text "102×4=" "835×2="
Replace-Text "576×4=" "319×8="
Replace-Text "975×6=" "378×7="
Replace-Text "843×5=" "194×8="
Replace-Text "877×3=" "156×5="
Replace-Text "206×3=" "937×6="
Replace-Text "846×6=" "288×6="
Replace-Text "443×7=" "909×8="
Replace-Text "937×2=" "574×3="
Replace-Text "104×8=" "689×4="
Replace-Text "423×9=" "655×9="
Replace-Text "709×4=" "429×5="
Replace-Text "565×8=" "793×9="
Replace-Text "428×3=" "394×6="
Replace-Text "860×7=" "935×5="
Replace-Text "522×6=" "234×5="
